# Weekly update: insert a new latest-week record at the top of the data
# (row 4) for "Agrícola del Norte S.A. de Arica - Ramas de apio", pushing
# all the older records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows (old row 4 .. old row 14) down by one row.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest week's record. Columns that stay
# the same as the previous top record (A, B, C, E, F, G, H, I, N, O, Q, R)
# are simply copied across; only the date and volume/price figures change.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44699
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112017
$ws.Cells.Item(4, 7).Value = "Ramas de apio"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 50
$ws.Cells.Item(4, 11).Value = 9000
$ws.Cells.Item(4, 12).Value = 9500
$ws.Cells.Item(4, 13).Value = 9250
$ws.Cells.Item(4, 14).Value = "`$/atado 7 kilos"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 9250
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Match the date style used by the other rows in column D.
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
